$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update outcome labels (title-cased / renamed values) in column B
$ws.Range("B2").Value = "Continued"
$ws.Range("B3").Value = "Transfer"
$ws.Range("B4").Value = "Return"
$ws.Range("B5").Value = "Death"
$ws.Range("B6").Value = "Other"
$ws.Range("B7").Value = "Previous"
$ws.Range("B8").Value = "Discharge"

# Update the active cell selection to B11 (was A11)
$ws.Range("B11").Select()
